# Updates the cryptos list on the active worksheet to the latest scraped values.
# Applies per-cell updates for the Coin / Link / Price / Volume(1h) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell='D2'; Value='29.508.90'},
    @{Cell='D3'; Value='1.852.63'},
    @{Cell='E4'; Value='  -0.10%  '},
    @{Cell='D5'; Value='243.03'},
    @{Cell='E5'; Value='  -1.21%  '},
    @{Cell='D6'; Value='0.6365'},
    @{Cell='E6'; Value='  -0.71%  '},
    @{Cell='D7'; Value='1.0000'},
    @{Cell='B8'; Value='OKB'},
    @{Cell='C8'; Value='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'},
    @{Cell='D8'; Value='48.30'},
    @{Cell='E8'; Value='  +2.02%  '},
    @{Cell='B9'; Value='Cardano'},
    @{Cell='C9'; Value='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'},
    @{Cell='D9'; Value='0.2987'},
    @{Cell='E9'; Value='  -0.52%  '},
    @{Cell='B10'; Value='Dogecoin'},
    @{Cell='C10'; Value='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'},
    @{Cell='D10'; Value='0.07467'},
    @{Cell='E10'; Value='  -0.70%  '},
    @{Cell='B11'; Value='Solana'},
    @{Cell='C11'; Value='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'},
    @{Cell='D11'; Value='24.26'},
    @{Cell='E11'; Value='  -0.58%  '},
    @{Cell='B12'; Value='TRON'},
    @{Cell='C12'; Value='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Cell='D12'; Value='0.07630'},
    @{Cell='E12'; Value='  -0.70%  '},
    @{Cell='B13'; Value='WrappedEther'},
    @{Cell='C13'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Cell='D13'; Value='1.849.35'},
    @{Cell='E13'; Value='  -1.13%  '},
    @{Cell='B14'; Value='Polkadot'},
    @{Cell='C14'; Value='https://coinranking.com/coin/25W7FG7om+polkadot-dot'},
    @{Cell='D14'; Value='5.027'},
    @{Cell='E14'; Value='  -1.00%  '},
    @{Cell='B15'; Value='Polygon'},
    @{Cell='C15'; Value='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'},
    @{Cell='D15'; Value='0.6849'},
    @{Cell='E15'; Value='  -0.68%  '},
    @{Cell='B16'; Value='Litecoin'},
    @{Cell='C16'; Value='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Cell='D16'; Value='83.55'},
    @{Cell='E16'; Value='  -1.09%  '},
    @{Cell='B17'; Value='ShibaInu'},
    @{Cell='C17'; Value='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'},
    @{Cell='D17'; Value='0.000009553'},
    @{Cell='E17'; Value='  +0.83%  '},
    @{Cell='B18'; Value='Uniswap'},
    @{Cell='C18'; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'},
    @{Cell='D18'; Value='6.147'},
    @{Cell='E18'; Value='  +0.76%  '},
    @{Cell='B19'; Value='WrappedBTC'},
    @{Cell='C19'; Value='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'},
    @{Cell='D19'; Value='29.505.42'},
    @{Cell='E19'; Value='  -1.38%  '},
    @{Cell='B20'; Value='WrappedliquidstakedEther2.0'},
    @{Cell='C20'; Value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{Cell='D20'; Value='2.070.94'},
    @{Cell='E20'; Value='  -2.77%  '},
    @{Cell='B21'; Value='BitcoinCash'},
    @{Cell='C21'; Value='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'},
    @{Cell='D21'; Value='234.67'},
    @{Cell='E21'; Value='  -2.71%  '},
    @{Cell='B22'; Value='Avalanche'},
    @{Cell='C22'; Value='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'},
    @{Cell='D22'; Value='12.54'},
    @{Cell='E22'; Value='  -1.36%  '},
    @{Cell='B23'; Value='Dai'},
    @{Cell='C23'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Cell='D23'; Value='0.9999'},
    @{Cell='E23'; Value='  -0.09%  '},
    @{Cell='B24'; Value='Chainlink'},
    @{Cell='C24'; Value='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'},
    @{Cell='D24'; Value='7.665'},
    @{Cell='E24'; Value='  +2.22%  '},
    @{Cell='B25'; Value='BinanceUSD'},
    @{Cell='C25'; Value='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'},
    @{Cell='D25'; Value='1.001'},
    @{Cell='E25'; Value='  -0.09%  '},
    @{Cell='B26'; Value='Monero'},
    @{Cell='C26'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Cell='D26'; Value='156.94'},
    @{Cell='E26'; Value='  -1.83%  '},
    @{Cell='B27'; Value='Stellar'},
    @{Cell='C27'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Cell='D27'; Value='0.1405'},
    @{Cell='B28'; Value='Cosmos'},
    @{Cell='C28'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell='D28'; Value='8.474'},
    @{Cell='E28'; Value='  -1.36%  '},
    @{Cell='B29'; Value='EthereumClassic'},
    @{Cell='C29'; Value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{Cell='D29'; Value='17.74'},
    @{Cell='E29'; Value='  -1.74%  '},
    @{Cell='B30'; Value='PancakeSwap'},
    @{Cell='C30'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{Cell='D30'; Value='1.483'},
    @{Cell='E30'; Value='  -1.40%  '},
    @{Cell='B31'; Value='Hedera'},
    @{Cell='C31'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D31'; Value='0.05994'},
    @{Cell='E31'; Value='  -2.48%  '},
    @{Cell='B32'; Value='Toncoin'},
    @{Cell='C32'; Value='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'},
    @{Cell='D32'; Value='1.250'},
    @{Cell='E32'; Value='  -2.00%  '},
    @{Cell='B33'; Value='Filecoin'},
    @{Cell='C33'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Cell='D33'; Value='4.121'},
    @{Cell='E33'; Value='  -1.32%  '},
    @{Cell='B34'; Value='InternetComputer(DFINITY)'},
    @{Cell='C34'; Value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'},
    @{Cell='D34'; Value='4.073'},
    @{Cell='E34'; Value='  -1.60%  '},
    @{Cell='D35'; Value='1.180'},
    @{Cell='E35'; Value='  +1.34%  '},
    @{Cell='B36'; Value='LidoDAOToken'},
    @{Cell='C36'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'},
    @{Cell='D36'; Value='1.861'},
    @{Cell='E36'; Value='  -0.55%  '},
    @{Cell='B37'; Value='ImmutableX'},
    @{Cell='C37'; Value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Cell='D37'; Value='0.7169'},
    @{Cell='E37'; Value='  -2.65%  '},
    @{Cell='B38'; Value='HuobiToken'},
    @{Cell='C38'; Value='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'},
    @{Cell='D38'; Value='2.598'},
    @{Cell='E38'; Value='  -0.34%  '},
    @{Cell='B39'; Value='MXToken'},
    @{Cell='C39'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell='D39'; Value='2.796'},
    @{Cell='E39'; Value='  -2.59%  '},
    @{Cell='B40'; Value='VeChain'},
    @{Cell='C40'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Cell='D40'; Value='0.01771'},
    @{Cell='E40'; Value='  -1.90%  '},
    @{Cell='B41'; Value='Maker'},
    @{Cell='C41'; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Cell='D41'; Value='1.198.15'},
    @{Cell='E41'; Value='  -2.15%  '},
    @{Cell='B42'; Value='TrustWalletToken'},
    @{Cell='C42'; Value='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'},
    @{Cell='D42'; Value='0.9092'},
    @{Cell='E42'; Value='  -2.27%  '},
    @{Cell='B43'; Value='FraxShare'},
    @{Cell='C43'; Value='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{Cell='D43'; Value='6.161'},
    @{Cell='E43'; Value='  -2.14%  '},
    @{Cell='B44'; Value='PaxDollar'},
    @{Cell='C44'; Value='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'},
    @{Cell='D44'; Value='0.9996'},
    @{Cell='E44'; Value='  -0.28%  '},
    @{Cell='B45'; Value='RocketPoolETH'},
    @{Cell='C45'; Value='https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'},
    @{Cell='D45'; Value='2.006.95'},
    @{Cell='E45'; Value='  -1.63%  '},
    @{Cell='B46'; Value='Quant'},
    @{Cell='C46'; Value='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'},
    @{Cell='D46'; Value='101.88'},
    @{Cell='E46'; Value='  -0.43%  '},
    @{Cell='B47'; Value='Aave'},
    @{Cell='C47'; Value='https://coinranking.com/coin/ixgUfzmLR+aave-aave'},
    @{Cell='D47'; Value='66.33'},
    @{Cell='E47'; Value='  -0.72%  '},
    @{Cell='D48'; Value='0.00000000123'},
    @{Cell='E48'; Value='  -0.40%  '},
    @{Cell='B49'; Value='Aptos'},
    @{Cell='C49'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Cell='D49'; Value='7.278'},
    @{Cell='E49'; Value='  +8.06%  '},
    @{Cell='B50'; Value='TheSandbox'},
    @{Cell='C50'; Value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Cell='D50'; Value='0.4031'},
    @{Cell='E50'; Value='  -1.81%  '},
    @{Cell='B51'; Value='EnergySwap'},
    @{Cell='C51'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='D51'; Value='9.103'},
    @{Cell='E51'; Value='  -3.05%  '}
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $col = $change.Cell.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "E") {
        # Force text storage so numeric-looking strings (e.g. "29.508.90",
        # "48.30", "1.0000") are kept verbatim instead of being parsed as
        # numbers, matching the original inline-string cell content.
        $cell.NumberFormat = "@"
        $cell.Value = $change.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $change.Value
    }
}

Write-Host "Applied $($changes.Count) cell updates"
